$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New annotation cells, written in the order that reproduces the
# --- author's shared-string table append order (TO DO: .. need to redesign...) ---

$ws.Cells.Item(18, 9).Value = "TO DO:"
$ws.Cells.Item(18, 9).WrapText = $true
$ws.Cells.Item(18, 9).Font.Bold = $true

$ws.Cells.Item(19, 9).Value = "Allow user to specify filename to load/save data."
$ws.Cells.Item(19, 9).WrapText = $true
$ws.Cells.Item(19, 9).Font.Bold = $true

$ws.Cells.Item(11, 10).Value = "with backup(s)"
$ws.Cells.Item(11, 10).WrapText = $true

$ws.Cells.Item(10, 10).Value = "clear lists"
$ws.Cells.Item(10, 10).WrapText = $true

$ws.Cells.Item(9, 10).Value = "Total value . Also by cases, trays, items, + materials."
$ws.Cells.Item(9, 10).WrapText = $true

$ws.Cells.Item(7, 10).Value = "not yet…"
$ws.Cells.Item(7, 10).WrapText = $true

$ws.Cells.Item(21, 9).Value = "Smart search"
$ws.Cells.Item(21, 9).WrapText = $true
$ws.Cells.Item(21, 9).Font.Bold = $true

$ws.Cells.Item(21, 10).Value = "need to redesign cases & trays to set item number and size limits. Should have done this at the beginning!"
$ws.Cells.Item(21, 10).WrapText = $true

# --- Row height for the wrapped, taller rows ---
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(21).RowHeight = 60

# --- Remove now-redundant checkbox cells ---
$ws.Cells.Item(21, 8).Clear()
$ws.Cells.Item(25, 8).Clear()
$ws.Cells.Item(26, 8).Clear()

# --- Flip some checkboxes from unchecked to checked ---
$ws.Cells.Item(22, 8).Value = $true
$ws.Cells.Item(23, 8).Value = $true
$ws.Cells.Item(24, 8).Value = $true

# --- "MyLinkedList" moves from row 28 down to row 29, restyled as a heading ---
$ws.Cells.Item(28, 6).Clear()
$ws.Cells.Item(28, 8).Clear()

$ws.Cells.Item(29, 6).Value = "MyLinkedList"
$ws.Cells.Item(29, 6).WrapText = $true
$ws.Cells.Item(29, 6).Font.Bold = $true
$ws.Cells.Item(29, 6).Font.Color = 255
$ws.Cells.Item(29, 8).Value = $true

# --- "JavaFX interface" restyled as a heading too ---
$ws.Cells.Item(30, 6).Font.Color = 255

# --- Sheet view: drop the pinned top-left cell, move the selection ---
$sheetView = $ws.Application.ActiveWindow
$ws.Range("A21").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# --- Workbook window geometry ---
$excel.ActiveWindow.WindowState = -4143
$excel.ActiveWindow.Left = 1290
$excel.ActiveWindow.Width = 27195
